$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $owner = $ws.Cells.Item($r, 3).Value()
    if ($owner -ne "Carol") {
        $ws.Cells.Item($r, 2).Value = "Critical"
        $ws.Cells.Item($r, 3).Value = "Carol"
    }
}
